$wb = $excel.ActiveWorkbook

$wsPerf = $wb.Worksheets.Item("Performance Metrics")
$wsPerf.Cells.Item(1,1).Value = 'Model Name'
$wsPerf.Cells.Item(1,2).Value = 'calinski_harabasz_score'
$wsPerf.Cells.Item(1,3).Value = 'silhouette_score'
$wsPerf.Cells.Item(1,4).Value = 'davies_bouldin_score'
$wsPerf.Cells.Item(1,5).Value = 'calinski_harabasz_score_Normalized'
$wsPerf.Cells.Item(1,6).Value = 'Model_Rank'

$wsPerf.Cells.Item(2,1).Value = 'agglomerativeClustering_SpectralEmbedding_BERT6'
$wsPerf.Cells.Item(2,2).Value = [double]'8.782187687781431e+30'
$wsPerf.Cells.Item(2,3).Value = [double]'0.5054226446462774'
$wsPerf.Cells.Item(2,4).Value = 0
$wsPerf.Cells.Item(2,5).Value = 1
$wsPerf.Cells.Item(2,6).Value = 1

$wsPerf.Cells.Item(3,1).Value = 'agglomerativeClustering_SpectralEmbedding_BERT7'
$wsPerf.Cells.Item(3,2).Value = [double]'8.351493156420156e+30'
$wsPerf.Cells.Item(3,3).Value = [double]'0.206780496129849'
$wsPerf.Cells.Item(3,4).Value = 0
$wsPerf.Cells.Item(3,5).Value = [double]'0.9509581727613843'
$wsPerf.Cells.Item(3,6).Value = 2

$wsPerf.Cells.Item(4,1).Value = 'agglomerativeClustering_SpectralEmbedding_BERT5'
$wsPerf.Cells.Item(4,2).Value = [double]'4.911648493278301e+30'
$wsPerf.Cells.Item(4,3).Value = [double]'0.5888765688896752'
$wsPerf.Cells.Item(4,4).Value = 0
$wsPerf.Cells.Item(4,5).Value = [double]'0.5592739153265681'
$wsPerf.Cells.Item(4,6).Value = 3

$wsClusters = $wb.Worksheets.Item("Clusters")
$wsClusters.Cells.Item(1,1).Value = 'Errors'
$wsClusters.Cells.Item(1,2).Value = 'agglomerativeClustering_SpectralEmbedding_BERT5'
$wsClusters.Cells.Item(1,3).Value = 'agglomerativeClustering_SpectralEmbedding_BERT6'
$wsClusters.Cells.Item(1,4).Value = 'agglomerativeClustering_SpectralEmbedding_BERT7'

$wsClusters.Cells.Item(2,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(2,2).Value = 0
$wsClusters.Cells.Item(2,3).Value = 1
$wsClusters.Cells.Item(2,4).Value = 0

$wsClusters.Cells.Item(3,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(3,2).Value = 2
$wsClusters.Cells.Item(3,3).Value = 2
$wsClusters.Cells.Item(3,4).Value = 2

$wsClusters.Cells.Item(4,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(4,2).Value = 4
$wsClusters.Cells.Item(4,3).Value = 4
$wsClusters.Cells.Item(4,4).Value = 4

$wsClusters.Cells.Item(5,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(5,2).Value = 0
$wsClusters.Cells.Item(5,3).Value = 1
$wsClusters.Cells.Item(5,4).Value = 0

$wsClusters.Cells.Item(6,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xa'']'
$wsClusters.Cells.Item(6,2).Value = 1
$wsClusters.Cells.Item(6,3).Value = 0
$wsClusters.Cells.Item(6,4).Value = 1

$wsClusters.Cells.Item(7,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xb'']'
$wsClusters.Cells.Item(7,2).Value = 1
$wsClusters.Cells.Item(7,3).Value = 0
$wsClusters.Cells.Item(7,4).Value = 1

$wsClusters.Cells.Item(8,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(8,2).Value = 4
$wsClusters.Cells.Item(8,3).Value = 4
$wsClusters.Cells.Item(8,4).Value = 4

$wsClusters.Cells.Item(9,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xb'']'
$wsClusters.Cells.Item(9,2).Value = 1
$wsClusters.Cells.Item(9,3).Value = 0
$wsClusters.Cells.Item(9,4).Value = 6

$wsClusters.Cells.Item(10,1).Value = '[''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x3d'', ''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x25'', ''LOCK Hard Hang, did not return to idle in 1s'']'
$wsClusters.Cells.Item(10,2).Value = 4
$wsClusters.Cells.Item(10,3).Value = 4
$wsClusters.Cells.Item(10,4).Value = 4

$wsClusters.Cells.Item(11,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(11,2).Value = 4
$wsClusters.Cells.Item(11,3).Value = 4
$wsClusters.Cells.Item(11,4).Value = 4

$wsClusters.Cells.Item(12,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xb'']'
$wsClusters.Cells.Item(12,2).Value = 1
$wsClusters.Cells.Item(12,3).Value = 0
$wsClusters.Cells.Item(12,4).Value = 1

$wsClusters.Cells.Item(13,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(13,2).Value = 4
$wsClusters.Cells.Item(13,3).Value = 4
$wsClusters.Cells.Item(13,4).Value = 4

$wsClusters.Cells.Item(14,1).Value = '[''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0xb'', ''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x5'', ''LOCK Hard Hang, did not return to idle in 1s'']'
$wsClusters.Cells.Item(14,2).Value = 4
$wsClusters.Cells.Item(14,3).Value = 4
$wsClusters.Cells.Item(14,4).Value = 4

$wsClusters.Cells.Item(15,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(15,2).Value = 0
$wsClusters.Cells.Item(15,3).Value = 1
$wsClusters.Cells.Item(15,4).Value = 0

$wsClusters.Cells.Item(16,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(16,2).Value = 0
$wsClusters.Cells.Item(16,3).Value = 1
$wsClusters.Cells.Item(16,4).Value = 0

$wsClusters.Cells.Item(17,1).Value = '[''mscod: MCE when CR4.MCE is clear-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''mscod: MCE under WFS-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x5'', ''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x49'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x3'', ''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(17,2).Value = 0
$wsClusters.Cells.Item(17,3).Value = 5
$wsClusters.Cells.Item(17,4).Value = 5

$wsClusters.Cells.Item(18,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(18,2).Value = 0
$wsClusters.Cells.Item(18,3).Value = 1
$wsClusters.Cells.Item(18,4).Value = 0

$wsClusters.Cells.Item(19,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x9'', ''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xb'']'
$wsClusters.Cells.Item(19,2).Value = 0
$wsClusters.Cells.Item(19,3).Value = 1
$wsClusters.Cells.Item(19,4).Value = 0

$wsClusters.Cells.Item(20,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xa'']'
$wsClusters.Cells.Item(20,2).Value = 1
$wsClusters.Cells.Item(20,3).Value = 0
$wsClusters.Cells.Item(20,4).Value = 1

$wsClusters.Cells.Item(21,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(21,2).Value = 0
$wsClusters.Cells.Item(21,3).Value = 1
$wsClusters.Cells.Item(21,4).Value = 0

$wsClusters.Cells.Item(22,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x9'', ''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x9'']'
$wsClusters.Cells.Item(22,2).Value = 0
$wsClusters.Cells.Item(22,3).Value = 1
$wsClusters.Cells.Item(22,4).Value = 0

$wsClusters.Cells.Item(23,1).Value = '[''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x11'', ''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x1b'', ''LOCK Hard Hang, did not return to idle in 1s'']'
$wsClusters.Cells.Item(23,2).Value = 4
$wsClusters.Cells.Item(23,3).Value = 4
$wsClusters.Cells.Item(23,4).Value = 4

$wsClusters.Cells.Item(24,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x9'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(24,2).Value = 0
$wsClusters.Cells.Item(24,3).Value = 1
$wsClusters.Cells.Item(24,4).Value = 0

$wsClusters.Cells.Item(25,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'']'
$wsClusters.Cells.Item(25,2).Value = 4
$wsClusters.Cells.Item(25,3).Value = 4
$wsClusters.Cells.Item(25,4).Value = 4

$wsClusters.Cells.Item(26,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x9'']'
$wsClusters.Cells.Item(26,2).Value = 1
$wsClusters.Cells.Item(26,3).Value = 0
$wsClusters.Cells.Item(26,4).Value = 1

$wsClusters.Cells.Item(27,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(27,2).Value = 3
$wsClusters.Cells.Item(27,3).Value = 3
$wsClusters.Cells.Item(27,4).Value = 3

$wsClusters.Cells.Item(28,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(28,2).Value = 3
$wsClusters.Cells.Item(28,3).Value = 3
$wsClusters.Cells.Item(28,4).Value = 3

$wsClusters.Cells.Item(29,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(29,2).Value = 3
$wsClusters.Cells.Item(29,3).Value = 3
$wsClusters.Cells.Item(29,4).Value = 3

$wsClusters.Cells.Item(30,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(30,2).Value = 3
$wsClusters.Cells.Item(30,3).Value = 3
$wsClusters.Cells.Item(30,4).Value = 3

$wsClusters.Cells.Item(31,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xb'']'
$wsClusters.Cells.Item(31,2).Value = 1
$wsClusters.Cells.Item(31,3).Value = 0
$wsClusters.Cells.Item(31,4).Value = 1

$wsClusters.Cells.Item(32,1).Value = '[''mscod: MCE when CR4.MCE is clear-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''mscod: MCE when CR4.MCE is clear-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xb'']'
$wsClusters.Cells.Item(32,2).Value = 4
$wsClusters.Cells.Item(32,3).Value = 4
$wsClusters.Cells.Item(32,4).Value = 4

$wsClusters.Cells.Item(33,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x5'', ''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xa'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x4'', ''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xb'']'
$wsClusters.Cells.Item(33,2).Value = 0
$wsClusters.Cells.Item(33,3).Value = 1
$wsClusters.Cells.Item(33,4).Value = 0

$wsClusters.Cells.Item(34,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(34,2).Value = 3
$wsClusters.Cells.Item(34,3).Value = 3
$wsClusters.Cells.Item(34,4).Value = 3

$wsClusters.Cells.Item(35,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(35,2).Value = 3
$wsClusters.Cells.Item(35,3).Value = 3
$wsClusters.Cells.Item(35,4).Value = 3

$wsClusters.Cells.Item(36,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xa'']'
$wsClusters.Cells.Item(36,2).Value = 0
$wsClusters.Cells.Item(36,3).Value = 1
$wsClusters.Cells.Item(36,4).Value = 0

$wsClusters.Cells.Item(37,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0xa'']'
$wsClusters.Cells.Item(37,2).Value = 1
$wsClusters.Cells.Item(37,3).Value = 0
$wsClusters.Cells.Item(37,4).Value = 6

$wsClusters.Cells.Item(38,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(38,2).Value = 3
$wsClusters.Cells.Item(38,3).Value = 3
$wsClusters.Cells.Item(38,4).Value = 3

$wsClusters.Cells.Item(39,1).Value = '[''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x1f'', ''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x17'', ''LOCK Hard Hang, did not return to idle in 1s'']'
$wsClusters.Cells.Item(39,2).Value = 4
$wsClusters.Cells.Item(39,3).Value = 4
$wsClusters.Cells.Item(39,4).Value = 4

$wsClusters.Cells.Item(40,1).Value = '[''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x9'', ''Merge Bank5-UPI signaled an MCA to Ubox; Check mc_status of MCA BANKID:5 and  MCA BANK_INDEX:0x1'', ''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0x1f'']'
$wsClusters.Cells.Item(40,2).Value = 0
$wsClusters.Cells.Item(40,3).Value = 1
$wsClusters.Cells.Item(40,4).Value = 0

$wsClusters.Cells.Item(41,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(41,2).Value = 4
$wsClusters.Cells.Item(41,3).Value = 4
$wsClusters.Cells.Item(41,4).Value = 4

$wsClusters.Cells.Item(42,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(42,2).Value = 1
$wsClusters.Cells.Item(42,3).Value = 0
$wsClusters.Cells.Item(42,4).Value = 6

$wsClusters.Cells.Item(43,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(43,2).Value = 1
$wsClusters.Cells.Item(43,3).Value = 0
$wsClusters.Cells.Item(43,4).Value = 6

$wsClusters.Cells.Item(44,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(44,2).Value = 1
$wsClusters.Cells.Item(44,3).Value = 0
$wsClusters.Cells.Item(44,4).Value = 6

$wsClusters.Cells.Item(45,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(45,2).Value = 4
$wsClusters.Cells.Item(45,3).Value = 4
$wsClusters.Cells.Item(45,4).Value = 4

$wsClusters.Cells.Item(46,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(46,2).Value = 1
$wsClusters.Cells.Item(46,3).Value = 0
$wsClusters.Cells.Item(46,4).Value = 6

$wsClusters.Cells.Item(47,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(47,2).Value = 1
$wsClusters.Cells.Item(47,3).Value = 0
$wsClusters.Cells.Item(47,4).Value = 6

$wsClusters.Cells.Item(48,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x0'']'
$wsClusters.Cells.Item(48,2).Value = 1
$wsClusters.Cells.Item(48,3).Value = 0
$wsClusters.Cells.Item(48,4).Value = 6

$wsClusters.Cells.Item(49,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'', ''Merge Bank7-CHA0 signaled an MCA to Ubox; Check mc_status of MCA BANKID:7 and  MCA BANK_INDEX:0xd'']'
$wsClusters.Cells.Item(49,2).Value = 4
$wsClusters.Cells.Item(49,3).Value = 4
$wsClusters.Cells.Item(49,4).Value = 4

$wsClusters.Cells.Item(50,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(50,2).Value = 2
$wsClusters.Cells.Item(50,3).Value = 2
$wsClusters.Cells.Item(50,4).Value = 2

$wsClusters.Cells.Item(51,1).Value = '[''Merge Bank6-Punit signaled an MCA to Ubox; Check mc_status of MCA BANKID:6 and  MCA BANK_INDEX:0x9'']'
$wsClusters.Cells.Item(51,2).Value = 1
$wsClusters.Cells.Item(51,3).Value = 0
$wsClusters.Cells.Item(51,4).Value = 6

$wsClusters.Cells.Item(52,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(52,2).Value = 2
$wsClusters.Cells.Item(52,3).Value = 2
$wsClusters.Cells.Item(52,4).Value = 2

$wsClusters.Cells.Item(53,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(53,2).Value = 2
$wsClusters.Cells.Item(53,3).Value = 2
$wsClusters.Cells.Item(53,4).Value = 2

$wsClusters.Cells.Item(54,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(54,2).Value = 2
$wsClusters.Cells.Item(54,3).Value = 2
$wsClusters.Cells.Item(54,4).Value = 2

$wsClusters.Cells.Item(55,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(55,2).Value = 2
$wsClusters.Cells.Item(55,3).Value = 2
$wsClusters.Cells.Item(55,4).Value = 2

$wsClusters.Cells.Item(56,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(56,2).Value = 2
$wsClusters.Cells.Item(56,3).Value = 2
$wsClusters.Cells.Item(56,4).Value = 2

$wsClusters.Cells.Item(57,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(57,2).Value = 2
$wsClusters.Cells.Item(57,3).Value = 2
$wsClusters.Cells.Item(57,4).Value = 2

$wsClusters.Cells.Item(58,1).Value = '[''mscod: MCE when MCIP bit is set-Error signaled by the core and logged in Ubox. Check the core for more details.'']'
$wsClusters.Cells.Item(58,2).Value = 2
$wsClusters.Cells.Item(58,3).Value = 2
$wsClusters.Cells.Item(58,4).Value = 2
